# feat: add 2022-Q3 data
#
# - "总计" sheet gains a new top data row for 2022-Q3 (existing 2022-Q2 / 2022-Q1
#   rows shift down and their running index is bumped).
# - A brand-new "2022-Q3" sheet is inserted right after "总计", holding the new
#   fund snapshot. The existing "2022-Q2" and "2022-Q1" sheets keep their own
#   name/content untouched, they simply shift one tab to the right.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) "总计" sheet: insert a fresh row 2 for 2022-Q3.
# ------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

$summary.Rows.Item(2).Insert()
$summary.Range("B2:D2").ClearFormats()

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 2
$summary.Range("D2").Value = 0.32

# Copy column-A's formatting (centered/bordered header style) down onto the
# freshly inserted cell, then refresh the running index on the rows that
# shifted down.
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)

$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2

# ------------------------------------------------------------------
# 2) New "2022-Q3" sheet, placed right after "总计".
#    Copy the existing "2022-Q2" sheet (item 2) so the new tab starts out
#    with identical formatting, then overwrite it with the Q3 fund data.
# ------------------------------------------------------------------
$oldQ2 = $wb.Worksheets.Item(2)
$oldQ2.Copy($oldQ2, $null)

$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

$q3.Range("B2:G3").NumberFormat = "@"

$q3.Range("B2").Value = "008980"
$q3.Range("C2").Value = "中邮科技创新精选混合A"
$q3.Range("D2").Value = "4.75"
$q3.Range("E2").Value = "88.40"
$q3.Range("F2").Value = "4.32"
$q3.Range("G2").Value = "0.2052"
$q3.Range("H2").Value = 10

$q3.Range("B3").Value = "008981"
$q3.Range("C3").Value = "中邮科技创新精选混合C"
$q3.Range("D3").Value = "2.66"
$q3.Range("E3").Value = "88.40"
$q3.Range("F3").Value = "4.32"
$q3.Range("G3").Value = "0.1149"
$q3.Range("H3").Value = 10

# Copying a sheet makes the copy the active one; restore the original
# active tab ("2022-Q1", now the 4th tab) so it stays the selected sheet,
# just like before the edit.
$wb.Worksheets.Item(4).Activate()
